# Update countries & provincias Spain
# Refresh COVID country stats (cases/recovered/deaths) and fix the
# "Datos actualizados" timestamp. A handful of rows also swap their
# country label because the source data got re-sorted by total cases
# (e.g. Moldavia/Serbia, Dinamarca/Bosnia, Siria/Jordania, Islas Feroe).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 13 de Agosto de 2020 a las 16:17"
# Row 4
$ws.Cells.Item(4, 2).Value = 5362253
$ws.Cells.Item(4, 3).Value = 1951
$ws.Cells.Item(4, 4).Value = 2813837
$ws.Cells.Item(4, 5).Value = 2379246
$ws.Cells.Item(4, 7).Value = 39
$ws.Cells.Item(4, 8).Value = 169170
# Row 6
$ws.Cells.Item(6, 2).Value = 2427066
$ws.Cells.Item(6, 3).Value = 31595
$ws.Cells.Item(6, 4).Value = 1721177
$ws.Cells.Item(6, 5).Value = 658417
$ws.Cells.Item(6, 7).Value = 334
$ws.Cells.Item(6, 8).Value = 47472
# Row 16
$ws.Cells.Item(16, 2).Value = 294519
$ws.Cells.Item(16, 3).Value = 1482
$ws.Cells.Item(16, 4).Value = 260393
$ws.Cells.Item(16, 5).Value = 30823
$ws.Cells.Item(16, 7).Value = 34
$ws.Cells.Item(16, 8).Value = 3303
# Row 19
$ws.Cells.Item(19, 4).Value = 192434
$ws.Cells.Item(19, 5).Value = 70894
$ws.Cells.Item(19, 7).Value = 33
$ws.Cells.Item(19, 8).Value = 5246
# Row 22
$ws.Cells.Item(22, 2).Value = 221413
$ws.Cells.Item(22, 3).Value = 563
$ws.Cells.Item(22, 5).Value = 11337
# Row 36
$ws.Cells.Item(36, 2).Value = 83852
$ws.Cells.Item(36, 8).Value = 5770
# Row 37
$ws.Cells.Item(37, 2).Value = 82531
$ws.Cells.Item(37, 3).Value = 232
$ws.Cells.Item(37, 4).Value = 77278
$ws.Cells.Item(37, 5).Value = 4702
$ws.Cells.Item(37, 7).Value = 12
$ws.Cells.Item(37, 8).Value = 551
# Row 45
$ws.Cells.Item(45, 2).Value = 61204
$ws.Cells.Item(45, 3).Value = 577
$ws.Cells.Item(45, 7).Value = 4
$ws.Cells.Item(45, 8).Value = 6165
# Row 47
$ws.Cells.Item(47, 4).Value = 50736
$ws.Cells.Item(47, 5).Value = 4734
# Row 54
$ws.Cells.Item(54, 2).Value = 41725
$ws.Cells.Item(54, 3).Value = 153
$ws.Cells.Item(54, 4).Value = 39495
$ws.Cells.Item(54, 5).Value = 2007
# Row 62
$ws.Cells.Item(62, 2).Value = 33093
$ws.Cells.Item(62, 3).Value = 628
$ws.Cells.Item(62, 4).Value = 26270
$ws.Cells.Item(62, 5).Value = 6608
$ws.Cells.Item(62, 7).Value = 5
$ws.Cells.Item(62, 8).Value = 215
# Row 64
$ws.Cells.Item(64, 1).Value = "Moldavia"
$ws.Cells.Item(64, 2).Value = 29087
$ws.Cells.Item(64, 3).Value = 390
$ws.Cells.Item(64, 4).Value = 19998
$ws.Cells.Item(64, 5).Value = 8211
$ws.Cells.Item(64, 7).Value = 15
$ws.Cells.Item(64, 8).Value = 878
# Row 65
$ws.Cells.Item(65, 1).Value = "Serbia"
$ws.Cells.Item(65, 2).Value = 28998
$ws.Cells.Item(65, 3).Value = 247
$ws.Cells.Item(65, 4).Value = 18965
$ws.Cells.Item(65, 5).Value = 9372
$ws.Cells.Item(65, 7).Value = 3
$ws.Cells.Item(65, 8).Value = 661
# Row 77
$ws.Cells.Item(77, 2).Value = 15491
$ws.Cells.Item(77, 3).Value = 307
$ws.Cells.Item(77, 4).Value = 9186
$ws.Cells.Item(77, 5).Value = 6200
# Row 78
$ws.Cells.Item(78, 1).Value = "Dinamarca"
$ws.Cells.Item(78, 2).Value = 15214
$ws.Cells.Item(78, 3).Value = 144
$ws.Cells.Item(78, 4).Value = 13131
$ws.Cells.Item(78, 5).Value = 1462
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 621
# Row 79
$ws.Cells.Item(79, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(79, 2).Value = 15184
$ws.Cells.Item(79, 3).Value = 223
$ws.Cells.Item(79, 4).Value = 9156
$ws.Cells.Item(79, 5).Value = 5570
$ws.Cells.Item(79, 7).Value = 5
$ws.Cells.Item(79, 8).Value = 458
# Row 83
$ws.Cells.Item(83, 2).Value = 12357
$ws.Cells.Item(83, 3).Value = 140
$ws.Cells.Item(83, 4).Value = 8662
$ws.Cells.Item(83, 5).Value = 3163
$ws.Cells.Item(83, 7).Value = 2
$ws.Cells.Item(83, 8).Value = 532
# Row 84
$ws.Cells.Item(84, 2).Value = 12115
$ws.Cells.Item(84, 3).Value = 82
$ws.Cells.Item(84, 4).Value = 6305
$ws.Cells.Item(84, 5).Value = 5018
$ws.Cells.Item(84, 7).Value = 6
$ws.Cells.Item(84, 8).Value = 792
# Row 86
$ws.Cells.Item(86, 2).Value = 9817
$ws.Cells.Item(86, 3).Value = 34
$ws.Cells.Item(86, 5).Value = 703
$ws.Cells.Item(86, 7).Value = 1
$ws.Cells.Item(86, 8).Value = 257
# Row 89
$ws.Cells.Item(89, 2).Value = 8663
$ws.Cells.Item(89, 3).Value = 162
$ws.Cells.Item(89, 4).Value = 7401
$ws.Cells.Item(89, 5).Value = 1016
# Row 94
$ws.Cells.Item(94, 2).Value = 7950
$ws.Cells.Item(94, 3).Value = 38
$ws.Cells.Item(94, 4).Value = 6741
$ws.Cells.Item(94, 5).Value = 1146
# Row 103
$ws.Cells.Item(103, 2).Value = 6050
$ws.Cells.Item(103, 3).Value = 180
$ws.Cells.Item(103, 4).Value = 5078
$ws.Cells.Item(103, 5).Value = 811
$ws.Cells.Item(103, 7).Value = 1
$ws.Cells.Item(103, 8).Value = 161
# Row 122
$ws.Cells.Item(122, 2).Value = 2882
$ws.Cells.Item(122, 3).Value = 1
$ws.Cells.Item(122, 5).Value = 225
# Row 134
$ws.Cells.Item(134, 2).Value = 1976
$ws.Cells.Item(134, 3).Value = 4
$ws.Cells.Item(134, 4).Value = 1852
$ws.Cells.Item(134, 5).Value = 114
# Row 141
$ws.Cells.Item(141, 1).Value = "Siria"
$ws.Cells.Item(141, 2).Value = 1402
$ws.Cells.Item(141, 3).Value = 75
$ws.Cells.Item(141, 4).Value = 395
$ws.Cells.Item(141, 5).Value = 954
$ws.Cells.Item(141, 8).Value = 53
# Row 142
$ws.Cells.Item(142, 1).Value = "Uruguay"
$ws.Cells.Item(142, 2).Value = 1393
$ws.Cells.Item(142, 3).Value = 0
$ws.Cells.Item(142, 4).Value = 1163
$ws.Cells.Item(142, 5).Value = 193
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 37
# Row 143
$ws.Cells.Item(143, 1).Value = "Uganda"
$ws.Cells.Item(143, 2).Value = 1353
$ws.Cells.Item(143, 3).Value = 21
$ws.Cells.Item(143, 4).Value = 1141
$ws.Cells.Item(143, 5).Value = 201
$ws.Cells.Item(143, 7).Value = 2
$ws.Cells.Item(143, 8).Value = 11
# Row 144
$ws.Cells.Item(144, 1).Value = "Jordania"
$ws.Cells.Item(144, 2).Value = 1320
$ws.Cells.Item(144, 3).Value = 17
$ws.Cells.Item(144, 4).Value = 1222
$ws.Cells.Item(144, 5).Value = 87
$ws.Cells.Item(144, 8).Value = 11
# Row 145
$ws.Cells.Item(145, 1).Value = "Letonia"
$ws.Cells.Item(145, 2).Value = 1307
$ws.Cells.Item(145, 3).Value = 4
$ws.Cells.Item(145, 4).Value = 1078
$ws.Cells.Item(145, 5).Value = 197
$ws.Cells.Item(145, 8).Value = 32
# Row 171
$ws.Cells.Item(171, 1).Value = "Islas Feroe"
$ws.Cells.Item(171, 2).Value = 362
$ws.Cells.Item(171, 3).Value = 23
$ws.Cells.Item(171, 4).Value = 225
$ws.Cells.Item(171, 5).Value = 137
$ws.Cells.Item(171, 8).Value = 0
# Row 172
$ws.Cells.Item(172, 1).Value = "Birmania"
$ws.Cells.Item(172, 2).Value = 361
$ws.Cells.Item(172, 4).Value = 318
$ws.Cells.Item(172, 5).Value = 37
$ws.Cells.Item(172, 8).Value = 6
# Row 173
$ws.Cells.Item(173, 1).Value = "Mauricio"
$ws.Cells.Item(173, 2).Value = 344
$ws.Cells.Item(173, 4).Value = 334
$ws.Cells.Item(173, 5).Value = 0
$ws.Cells.Item(173, 8).Value = 10
# Row 184
$ws.Cells.Item(184, 2).Value = 205
$ws.Cells.Item(184, 3).Value = 2
$ws.Cells.Item(184, 4).Value = 188
$ws.Cells.Item(184, 5).Value = 17
